# Auto-generated edit script: updates calculated market-profit columns (H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW based on refreshed price data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3746.0476  # H40 was 3666.682
$ws.Cells.Item(40, 10).Value = 4183.1665  # J40 was 4015.2307
$ws.Cells.Item(40, 12).Value = 4183.1665  # L40 was 4015.2307
$ws.Cells.Item(40, 14).Value = -4533.1665  # N40 was -4365.2307
$ws.Cells.Item(70, 8).Value = 5216.6177  # H70 was 5323.485
$ws.Cells.Item(70, 9).Value = 8021.923  # I70 was 8022.385
$ws.Cells.Item(70, 10).Value = 3480  # J70 was 3569.2
$ws.Cells.Item(70, 11).Value = 24065.769  # K70 was 24067.155
$ws.Cells.Item(70, 12).Value = 10440  # L70 was 10707.6
$ws.Cells.Item(70, 13).Value = -23795.769  # M70 was -23797.155
$ws.Cells.Item(70, 14).Value = -10980  # N70 was -11247.6
$ws.Cells.Item(73, 8).Value = 5216.6177  # H73 was 5323.485
$ws.Cells.Item(73, 9).Value = 8021.923  # I73 was 8022.385
$ws.Cells.Item(73, 10).Value = 3480  # J73 was 3569.2
$ws.Cells.Item(73, 11).Value = 24065.769  # K73 was 24067.155
$ws.Cells.Item(73, 12).Value = 10440  # L73 was 10707.6
$ws.Cells.Item(73, 13).Value = -23129.769  # M73 was -23131.155
$ws.Cells.Item(73, 14).Value = -12312  # N73 was -12579.6
$ws.Cells.Item(80, 8).Value = 17596.2  # H80 was 18836.822
$ws.Cells.Item(80, 9).Value = 6125.6665  # I80 was 6475.9414
$ws.Cells.Item(80, 10).Value = 34802  # J80 was 37940
$ws.Cells.Item(80, 11).Value = 18376.9995  # K80 was 19427.8242
$ws.Cells.Item(80, 12).Value = 104406  # L80 was 113820
$ws.Cells.Item(80, 13).Value = -17378.9995  # M80 was -18429.8242
$ws.Cells.Item(80, 14).Value = -106402  # N80 was -115816
$ws.Cells.Item(83, 8).Value = 17596.2  # H83 was 18836.822
$ws.Cells.Item(83, 9).Value = 6125.6665  # I83 was 6475.9414
$ws.Cells.Item(83, 10).Value = 34802  # J83 was 37940
$ws.Cells.Item(83, 11).Value = 55130.9985  # K83 was 58283.47259999999
$ws.Cells.Item(83, 12).Value = 313218  # L83 was 341460
$ws.Cells.Item(83, 13).Value = -50138.9985  # M83 was -53291.47259999999
$ws.Cells.Item(83, 14).Value = -323202  # N83 was -351444
$ws.Cells.Item(115, 8).Value = 669.2727  # H115 was 669.8182
$ws.Cells.Item(115, 9).Value = 636.3  # I115 was 636.9
$ws.Cells.Item(115, 11).Value = 1908.9  # K115 was 1910.7
$ws.Cells.Item(115, 13).Value = -341.8999999999999  # M115 was -343.6999999999998
$ws.Cells.Item(116, 8).Value = 10422204  # H116 was 10005439
$ws.Cells.Item(116, 9).Value = 19234528  # I116 was 17860852
$ws.Cells.Item(116, 11).Value = 19234528  # K116 was 17860852
$ws.Cells.Item(116, 13).Value = -19231086  # M116 was -17857410
$ws.Cells.Item(125, 8).Value = 62501224  # H125 was 83334790
$ws.Cells.Item(125, 9).Value = 250000020  # I125 was 500000000
$ws.Cells.Item(125, 10).Value = 1625  # J125 was 1750
$ws.Cells.Item(125, 11).Value = 2250000180  # K125 was 4500000000
$ws.Cells.Item(125, 12).Value = 14625  # L125 was 15750
$ws.Cells.Item(125, 13).Value = -2249997720  # M125 was -4499997540
$ws.Cells.Item(125, 14).Value = -19545  # N125 was -20670
$ws.Cells.Item(132, 8).Value = 1165.3049  # H132 was 1151.6786
$ws.Cells.Item(132, 9).Value = 1141.3125  # I132 was 1136.8518
$ws.Cells.Item(132, 10).Value = 2125  # J132 was 1552
$ws.Cells.Item(132, 11).Value = 3423.9375  # K132 was 3410.5554
$ws.Cells.Item(132, 12).Value = 6375  # L132 was 4656
$ws.Cells.Item(132, 13).Value = -893.9375  # M132 was -880.5553999999997
$ws.Cells.Item(132, 14).Value = -11435  # N132 was -9716
$ws.Cells.Item(133, 8).Value = 12612492  # H133 was 8441922
$ws.Cells.Item(133, 10).Value = 12612492  # J133 was 8441922
$ws.Cells.Item(133, 12).Value = 12612492  # L133 was 8441922
$ws.Cells.Item(133, 14).Value = -12622612  # N133 was -8452042
$ws.Cells.Item(135, 8).Value = 2001379.8  # H135 was 2001390
$ws.Cells.Item(135, 9).Value = 2501474.8  # I135 was 3334983.2
$ws.Cells.Item(135, 11).Value = 22513273.2  # K135 was 30014848.8
$ws.Cells.Item(135, 13).Value = -22510738.2  # M135 was -30012313.8
$ws.Cells.Item(141, 8).Value = 4198.8  # H141 was 3132.8333
$ws.Cells.Item(141, 9).Value = 4198.8  # I141 was 3866.6667
$ws.Cells.Item(141, 10).Value = 0  # J141 was 2399
$ws.Cells.Item(141, 11).Value = 12596.4  # K141 was 11600.0001
$ws.Cells.Item(141, 12).Value = 0  # L141 was 7197
$ws.Cells.Item(141, 13).Value = -7416.400000000001  # M141 was -6420.000100000001
$ws.Cells.Item(141, 14).ClearContents()  # N141 was -17557

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 52633590  # H2 was 45456350
$ws.Cells.Item(2, 9).Value = 1686.2142  # I2 was 1473.2354
$ws.Cells.Item(2, 11).Value = 1686.2142  # K2 was 1473.2354
$ws.Cells.Item(2, 13).Value = -1573.2142  # M2 was -1360.2354
$ws.Cells.Item(61, 8).Value = 9116.0625  # H61 was 8625.764999999999
$ws.Cells.Item(61, 9).Value = 1406  # I61 was 1327.875
$ws.Cells.Item(61, 11).Value = 1406  # K61 was 1327.875
$ws.Cells.Item(61, 13).Value = -1194  # M61 was -1115.875
$ws.Cells.Item(102, 8).Value = 1901.9  # H102 was 2023.2222
$ws.Cells.Item(102, 9).Value = 1901.9  # I102 was 2023.2222
$ws.Cells.Item(102, 11).Value = 1901.9  # K102 was 2023.2222
$ws.Cells.Item(102, 13).Value = -279.9000000000001  # M102 was -401.2221999999999
$ws.Cells.Item(116, 8).Value = 52633590  # H116 was 45456350
$ws.Cells.Item(116, 9).Value = 1686.2142  # I116 was 1473.2354
$ws.Cells.Item(116, 11).Value = 1686.2142  # K116 was 1473.2354
$ws.Cells.Item(116, 13).Value = 607.7858000000001  # M116 was 820.7646
$ws.Cells.Item(132, 8).Value = 5718.7173  # H132 was 5803.8223
$ws.Cells.Item(132, 9).Value = 3482.739  # I132 was 3555.182
$ws.Cells.Item(132, 11).Value = 10448.217  # K132 was 10665.546
$ws.Cells.Item(132, 13).Value = -7918.217000000001  # M132 was -8135.545999999998
$ws.Cells.Item(136, 8).Value = 9116.0625  # H136 was 8625.764999999999
$ws.Cells.Item(136, 9).Value = 1406  # I136 was 1327.875
$ws.Cells.Item(136, 11).Value = 4218  # K136 was 3983.625
$ws.Cells.Item(136, 13).Value = -1668  # M136 was -1433.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 52633590  # H3 was 45456350
$ws.Cells.Item(3, 9).Value = 1686.2142  # I3 was 1473.2354
$ws.Cells.Item(3, 11).Value = 1686.2142  # K3 was 1473.2354
$ws.Cells.Item(3, 13).Value = -1572.2142  # M3 was -1359.2354
$ws.Cells.Item(8, 8).Value = 5000  # H8 was 9900
$ws.Cells.Item(8, 9).Value = 5000  # I8 was 9900
$ws.Cells.Item(8, 11).Value = 5000  # K8 was 9900
$ws.Cells.Item(8, 13).Value = -4860  # M8 was -9760
$ws.Cells.Item(134, 8).Value = 7511.6787  # H134 was 7719.6294
$ws.Cells.Item(134, 9).Value = 2379.1  # I134 was 2432.6667
$ws.Cells.Item(134, 11).Value = 7137.299999999999  # K134 was 7298.000100000001
$ws.Cells.Item(134, 13).Value = -4602.299999999999  # M134 was -4763.000100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 5275.7646  # H16 was 5528.375
$ws.Cells.Item(16, 9).Value = 1706  # I16 was 1824
$ws.Cells.Item(16, 11).Value = 1706  # K16 was 1824
$ws.Cells.Item(16, 13).Value = -1419  # M16 was -1537
$ws.Cells.Item(31, 8).Value = 7579.25  # H31 was 7747.6665
$ws.Cells.Item(31, 9).Value = 2466.3333  # I31 was 2551.9412
$ws.Cells.Item(31, 11).Value = 2466.3333  # K31 was 2551.9412
$ws.Cells.Item(31, 13).Value = -2171.3333  # M31 was -2256.9412
$ws.Cells.Item(34, 8).Value = 7579.25  # H34 was 7747.6665
$ws.Cells.Item(34, 9).Value = 2466.3333  # I34 was 2551.9412
$ws.Cells.Item(34, 11).Value = 2466.3333  # K34 was 2551.9412
$ws.Cells.Item(34, 13).Value = -2264.3333  # M34 was -2349.9412
$ws.Cells.Item(58, 8).Value = 7889.8887  # H58 was 8458.6
$ws.Cells.Item(58, 9).Value = 1890.8889  # I58 was 2208
$ws.Cells.Item(58, 11).Value = 1890.8889  # K58 was 2208
$ws.Cells.Item(58, 13).Value = -1687.8889  # M58 was -2005
$ws.Cells.Item(94, 8).Value = 1483.2307  # H94 was 1405.7858
$ws.Cells.Item(94, 10).Value = 1061.1111  # J94 was 994.9
$ws.Cells.Item(94, 12).Value = 1061.1111  # L94 was 994.9
$ws.Cells.Item(94, 14).Value = -1963.1111  # N94 was -1896.9
$ws.Cells.Item(99, 8).Value = 3749.6667  # H99 was 4068.5652
$ws.Cells.Item(99, 9).Value = 2772.1365  # I99 was 3016.9412
$ws.Cells.Item(99, 10).Value = 6437.875  # J99 was 7048.1665
$ws.Cells.Item(99, 11).Value = 2772.1365  # K99 was 3016.9412
$ws.Cells.Item(99, 12).Value = 6437.875  # L99 was 7048.1665
$ws.Cells.Item(99, 13).Value = -1274.1365  # M99 was -1518.9412
$ws.Cells.Item(99, 14).Value = -9433.875  # N99 was -10044.1665
$ws.Cells.Item(113, 8).Value = 5275.7646  # H113 was 5528.375
$ws.Cells.Item(113, 9).Value = 1706  # I113 was 1824
$ws.Cells.Item(113, 11).Value = 1706  # K113 was 1824
$ws.Cells.Item(113, 13).Value = 464  # M113 was 346
$ws.Cells.Item(126, 8).Value = 3749.6667  # H126 was 4068.5652
$ws.Cells.Item(126, 9).Value = 2772.1365  # I126 was 3016.9412
$ws.Cells.Item(126, 10).Value = 6437.875  # J126 was 7048.1665
$ws.Cells.Item(126, 11).Value = 8316.4095  # K126 was 9050.8236
$ws.Cells.Item(126, 12).Value = 19313.625  # L126 was 21144.4995
$ws.Cells.Item(126, 13).Value = -5846.4095  # M126 was -6580.8236
$ws.Cells.Item(126, 14).Value = -24253.625  # N126 was -26084.4995
$ws.Cells.Item(134, 8).Value = 7449.033  # H134 was 7465.8667
$ws.Cells.Item(134, 9).Value = 3333.4285  # I134 was 3430.8333
$ws.Cells.Item(134, 10).Value = 8701.608  # J134 was 8474.625
$ws.Cells.Item(134, 11).Value = 10000.2855  # K134 was 10292.4999
$ws.Cells.Item(134, 12).Value = 26104.824  # L134 was 25423.875
$ws.Cells.Item(134, 13).Value = -7465.2855  # M134 was -7757.499899999999
$ws.Cells.Item(134, 14).Value = -31174.824  # N134 was -30493.875
$ws.Cells.Item(136, 8).Value = 7889.8887  # H136 was 8458.6
$ws.Cells.Item(136, 9).Value = 1890.8889  # I136 was 2208
$ws.Cells.Item(136, 11).Value = 5672.6667  # K136 was 6624
$ws.Cells.Item(136, 13).Value = -3122.6667  # M136 was -4074

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 113952.18  # H2 was 113955.59
$ws.Cells.Item(2, 10).Value = 275668.25  # J2 was 275677.62
$ws.Cells.Item(2, 12).Value = 1654009.5  # L2 was 1654065.72
$ws.Cells.Item(2, 14).Value = -1654235.5  # N2 was -1654291.72
$ws.Cells.Item(4, 8).Value = 4910087.5  # H4 was 4686916
$ws.Cells.Item(4, 10).Value = 297.55554  # J4 was 299.7
$ws.Cells.Item(4, 12).Value = 892.66662  # L4 was 899.0999999999999
$ws.Cells.Item(4, 14).Value = -1116.66662  # N4 was -1123.1
$ws.Cells.Item(121, 8).Value = 1322.421  # H121 was 1640.5385
$ws.Cells.Item(121, 9).Value = 1399.4546  # I121 was 1513.1
$ws.Cells.Item(121, 10).Value = 1216.5  # J121 was 2065.3333
$ws.Cells.Item(121, 11).Value = 4198.3638  # K121 was 4539.299999999999
$ws.Cells.Item(121, 12).Value = 3649.5  # L121 was 6195.999899999999
$ws.Cells.Item(121, 13).Value = -2888.3638  # M121 was -3229.299999999999
$ws.Cells.Item(121, 14).Value = -6269.5  # N121 was -8815.999899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 667650.5  # H107 was 667650.9399999999
$ws.Cells.Item(107, 9).Value = 728073.25  # I107 was 728073.75
$ws.Cells.Item(107, 11).Value = 728073.25  # K107 was 728073.75
$ws.Cells.Item(107, 13).Value = -726153.25  # M107 was -726153.75
$ws.Cells.Item(122, 8).Value = 3827051  # H122 was 4276717.5
$ws.Cells.Item(122, 9).Value = 14526463  # I122 was 18157642
$ws.Cells.Item(122, 10).Value = 5832.2856  # J122 was 5664.077
$ws.Cells.Item(122, 11).Value = 43579389  # K122 was 54472926
$ws.Cells.Item(122, 12).Value = 17496.8568  # L122 was 16992.231
$ws.Cells.Item(122, 13).Value = -43576939  # M122 was -54470476
$ws.Cells.Item(122, 14).Value = -22396.8568  # N122 was -21892.231

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value = 44110  # H94 was 43582.5
$ws.Cells.Item(94, 10).Value = 44110  # J94 was 43582.5
$ws.Cells.Item(94, 12).Value = 44110  # L94 was 43582.5
$ws.Cells.Item(94, 14).Value = -45462  # N94 was -44934.5
$ws.Cells.Item(136, 8).Value = 15453  # H136 was 15962.667
$ws.Cells.Item(136, 9).Value = 3644  # I136 was 3737.8572
$ws.Cells.Item(136, 10).Value = 19704.24  # J136 was 19253.96
$ws.Cells.Item(136, 11).Value = 10932  # K136 was 11213.5716
$ws.Cells.Item(136, 12).Value = 59112.72  # L136 was 57761.88
$ws.Cells.Item(136, 13).Value = -8382  # M136 was -8663.571599999999
$ws.Cells.Item(136, 14).Value = -64212.72  # N136 was -62861.88
